$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-7 (columns E-H and M-T)
$data = @{
    2 = @{ E=2; F=1; G=1.6164125; H=3.232825; M=16.558025; N=33.11605; O=0.2047872600336892; P=0.155585835247525; Q=26.7645985853125; R=107.05839434125; S=0.2047872600336892; T=0.155585835247525 }
    3 = @{ E=2; F=1; G=1.6164125; H=3.232825; M=20.31779433333333; N=60.953383; O=0.2512875437409564; P=0.2863712008291233; Q=32.84193673282916; R=197.051620396975; S=0.2512875437409564; T=0.2863712008291233 }
    4 = @{ E=2; F=1; G=1.6164125; H=3.232825; M=10.405091; N=31.215273; O=0.1286886616182304; P=0.1466556042216543; Q=16.8189191560375; R=100.913514936225; S=0.1286886616182304; T=0.1466556042216543 }
    5 = @{ E=2; F=1; G=1.6164125; H=3.232825; M=13.1587975; N=26.317595; O=0.1627461056112162; P=0.1236453320906656; Q=21.27004476396875; R=85.08017905587499; S=0.1627461056112162; T=0.1236453320906656 }
    6 = @{ E=2; F=1; G=1.6164125; H=3.232825; M=7.303315666666667; N=21.909947; O=0.09032635260170116; P=0.1029373190408882; Q=11.80517073504583; R=70.83102441027501; S=0.09032635260170116; T=0.1029373190408882 }
    7 = @{ E=2; F=1; G=1.6164125; H=3.232825; M=13.11173766666667; N=39.335213; O=0.1621640763942067; P=0.1848047085701436; Q=21.19397666112084; R=127.163859966725; S=0.1621640763942067; T=0.1848047085701436 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
